$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, whether it needs the
# 'force text' guard (value looks like a plain number and Excel
# would otherwise auto-convert it from text to a numeric type).
$changes = @(
    @('D2', '27.770.34', $false),
    @('E2', '  -0.29%  ', $false),
    @('D3', '1.595.84', $false),
    @('E3', '  -1.55%  ', $false),
    @('E4', '  +0.04%  ', $false),
    @('D5', '209.16', $true),
    @('E5', '  -0.98%  ', $false),
    @('E7', '  +0.04%  ', $false),
    @('D8', '22.37', $true),
    @('E8', '  -2.65%  ', $false),
    @('E9', '  -1.40%  ', $false),
    @('D10', '0.0594', $true),
    @('E10', '  -1.59%  ', $false),
    @('E11', '  -1.40%  ', $false),
    @('D12', '1.822.87', $false),
    @('E12', '  -1.64%  ', $false),
    @('D13', '1.592.06', $false),
    @('E13', '  -1.83%  ', $false),
    @('E14', '  -2.19%  ', $false),
    @('E15', '  -3.22%  ', $false),
    @('D16', '27.759.57', $false),
    @('E16', '  -0.38%  ', $false),
    @('D17', '63.49', $true),
    @('E17', '  -1.48%  ', $false),
    @('D18', '219.72', $true),
    @('E18', '  -2.80%  ', $false),
    @('E19', '  -2.03%  ', $false),
    @('E20', '  -2.49%  ', $false),
    @('E21', '  +0.02%  ', $false),
    @('D22', '4.16', $true),
    @('E22', '  -3.34%  ', $false),
    @('D23', '9.81', $true),
    @('E23', '  -1.13%  ', $false),
    @('D24', '1.98', $true),
    @('E24', '  -3.89%  ', $false),
    @('D25', '154.06', $true),
    @('E25', '  +0.00%  ', $false),
    @('D26', '7.22', $true),
    @('E26', '  +4.86%  ', $false),
    @('E27', '  +0.05%  ', $false),
    @('D28', '15.18', $true),
    @('E28', '  -0.76%  ', $false),
    @('D29', '0.105', $true),
    @('E29', '  -3.65%  ', $false),
    @('D30', '1.16', $true),
    @('E30', '  -1.26%  ', $false),
    @('E31', '  -1.12%  ', $false),
    @('E32', '  -3.87%  ', $false),
    @('D33', '1.376.75', $false),
    @('E33', '  -2.86%  ', $false),
    @('E34', '  -2.73%  ', $false),
    @('E35', '  -3.20%  ', $false),
    @('D36', '0.974', $true),
    @('E36', '  +0.13%  ', $false),
    @('E37', '  +0.12%  ', $false),
    @('E38', '  -0.29%  ', $false),
    @('E39', '  -2.64%  ', $false),
    @('E40', '  -1.71%  ', $false),
    @('E41', '  +0.13%  ', $false),
    @('D42', '0.972', $true),
    @('E42', '  -2.76%  ', $false),
    @('D43', '64.61', $true),
    @('E43', '  -0.52%  ', $false),
    @('E44', '  +2.51%  ', $false),
    @('B45', 'FraxShare', $false),
    @('C45', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', $false),
    @('D45', '5.26', $true),
    @('E45', '  -1.59%  ', $false),
    @('B46', 'RenderToken', $false),
    @('C46', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', $false),
    @('D46', '1.75', $true),
    @('E46', '  -1.35%  ', $false),
    @('D47', '1.734.18', $false),
    @('E47', '  -1.66%  ', $false),
    @('D48', '86.67', $true),
    @('E48', '  -3.01%  ', $false),
    @('E49', '  -0.50%  ', $false),
    @('D50', '0.0968', $true),
    @('E50', '  -2.74%  ', $false),
    @('E51', '  -0.76%  ', $false)
)

foreach ($item in $changes) {
    $addr = $item[0]
    $val = $item[1]
    $needsTextGuard = $item[2]
    $cell = $ws.Range($addr)
    if ($needsTextGuard) {
        # Temporarily force a text number-format so Excel stores the
        # value as a string instead of silently parsing it into a
        # double, then restore the default style so no stray
        # formatting is left behind on the cell.
        $cell.NumberFormat = '@'
        $cell.Value = $val
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $val
    }
}

Write-Host ("Applied " + $changes.Count + " cell updates")
